$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Tue Nov 12 17:51:47 EST 2024"
$ws.Range("B3").Value = "Tue Nov 12 17:51:59 EST 2024"
$ws.Range("B4").Value = "Tue Nov 12 17:52:11 EST 2024"
$ws.Range("B5").Value = "Tue Nov 12 17:52:25 EST 2024"
$ws.Range("B6").Value = "Tue Nov 12 17:52:37 EST 2024"
$ws.Range("B7").Value = "Tue Nov 12 17:52:49 EST 2024"
$ws.Range("B8").Value = "Tue Nov 12 17:53:01 EST 2024"
$ws.Range("B9").Value = "Tue Nov 12 17:53:13 EST 2024"
$ws.Range("B10").Value = "Tue Nov 12 17:53:24 EST 2024"
$ws.Range("B11").Value = "Tue Nov 12 17:53:36 EST 2024"
$ws.Range("B12").Value = "Tue Nov 12 17:53:50 EST 2024"
$ws.Range("B13").Value = "Tue Nov 12 17:54:02 EST 2024"
$ws.Range("B14").Value = "Tue Nov 12 17:54:15 EST 2024"
$ws.Range("B15").Value = "Tue Nov 12 17:54:27 EST 2024"
$ws.Range("B16").Value = "Tue Nov 12 17:54:41 EST 2024"
$ws.Range("B17").Value = "Tue Nov 12 17:54:53 EST 2024"
$ws.Range("B18").Value = "Tue Nov 12 17:55:06 EST 2024"
$ws.Range("B19").Value = "Tue Nov 12 17:55:19 EST 2024"
$ws.Range("B20").Value = "Tue Nov 12 17:55:33 EST 2024"
$ws.Range("B21").Value = "Tue Nov 12 17:55:49 EST 2024"
$ws.Range("B22").Value = "Tue Nov 12 17:56:02 EST 2024"
$ws.Range("B23").Value = "Tue Nov 12 17:56:15 EST 2024"
$ws.Range("B24").Value = "Tue Nov 12 17:56:28 EST 2024"
$ws.Range("B25").Value = "Tue Nov 12 17:56:42 EST 2024"
$ws.Range("B26").Value = "Tue Nov 12 17:56:55 EST 2024"
$ws.Range("B27").Value = "Tue Nov 12 17:57:08 EST 2024"
$ws.Range("B28").Value = "Tue Nov 12 17:57:21 EST 2024"
$ws.Range("B29").Value = "Tue Nov 12 17:57:37 EST 2024"
$ws.Range("B30").Value = "Tue Nov 12 17:57:52 EST 2024"
$ws.Range("B31").Value = "Tue Nov 12 17:58:05 EST 2024"
$ws.Range("B32").Value = "Tue Nov 12 17:58:18 EST 2024"
$ws.Range("B33").Value = "Tue Nov 12 17:58:31 EST 2024"
$ws.Range("B34").Value = "Tue Nov 12 17:58:44 EST 2024"
$ws.Range("B36").Value = "Tue Nov 12 17:58:57 EST 2024"
$ws.Range("B37").Value = "Tue Nov 12 17:59:09 EST 2024"
$ws.Range("B38").Value = "Tue Nov 12 17:59:20 EST 2024"
$ws.Range("B39").Value = "Tue Nov 12 17:59:32 EST 2024"
$ws.Range("B40").Value = "Tue Nov 12 17:59:46 EST 2024"
$ws.Range("B41").Value = "Tue Nov 12 17:59:58 EST 2024"
$ws.Range("B42").Value = "Tue Nov 12 18:00:10 EST 2024"
$ws.Range("B43").Value = "Tue Nov 12 18:00:23 EST 2024"
$ws.Range("B44").Value = "Tue Nov 12 18:00:37 EST 2024"
$ws.Range("B45").Value = "Tue Nov 12 18:00:52 EST 2024"
$ws.Range("B47").Value = "Tue Nov 12 18:01:06 EST 2024"
$ws.Range("B48").Value = "Tue Nov 12 18:01:19 EST 2024"
$ws.Range("B49").Value = "Tue Nov 12 18:01:32 EST 2024"
$ws.Range("B50").Value = "Tue Nov 12 18:01:48 EST 2024"
$ws.Range("B51").Value = "Tue Nov 12 18:02:01 EST 2024"
$ws.Range("B52").Value = "Tue Nov 12 18:02:14 EST 2024"
$ws.Range("B53").Value = "Tue Nov 12 18:02:31 EST 2024"
$ws.Range("B54").Value = "Tue Nov 12 18:02:45 EST 2024"
